{"js": "// Insert the new \"leta\u010dki dan\" (flight day) log entries for days 11-13\n// right after the existing \"Deseti leta\u010dki dan 6.4.2022.\" paragraph.\nconst newParagraphs = [\n  {\n    \"text\": \"U ovom danu nastavljalo se tunirati regulatore prema navedenim uputama. Regulatori su se tunirali u Postion modeu. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U prvom letu (log_0) VTOL letjelica je u\u0161la u nestabilno pona\u0161anje zbog promjene parametara regulatora prije po\u010detka leta. Zaklju\u010deno je da se to dogodilo zbog promjene poja\u010danja na proporcionalnom dijelu pitch rate regulatora te je on nakon leta vra\u0107en na prija\u0161nju vrijednost. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U drugom letu (log_2) provjeralo se pona\u0161anje letjelice sa staarim parametrima regulatora te je let na kraju bio stabilan \u0161to se i o\u010dekivalo. U ovo letu prou\u010deni su utjecaji promjene parametara vezani uz orijentaciju plo\u010dice, no zaklju\u010deno je kako se ti parametri automatski odrede prilikom kalibracije pa ne postoji potreba za njihovim mijenjanjem. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U tre\u0107em letu (log_5) nastavilo se s tuniranjem regulatora. Tunirali su se parametri vezani uz pitch i pitchrate odnosno: MC_PITCHRATE_K, MC_PITCH_RATE_D, MC_PITCHRATE_I te MC_PITCH_P. Oni su se mijenjali sve dok se nije dobilo zadovoljavaju\u0107e pona\u0161anje po pitch odnosno dok letjelica nije morala sletiti zbog postotka baterije. Tako\u0111er ovaj dan su spremljeni parametri i stavljeni u share folder gdje se nalaze i logovi. Ovi parametri su karakterizirani kao 'OK parametri' te su se koristili u slu\u010daju kad bi se prilikom tuniranja do\u0161lo do kombinacije parametara pri kojoj bi pona\u0161anje bilo lo\u0161e ili nestabilno.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"Jedanaesti leta\u010dki dan 7.4.2022.\",\n    \"list\": true\n  },\n  {\n    \"text\": \"Nastavilo se sa tuniranjem regulatora kako bi se dobili jo\u0161 bolji parametri odnosno kako bi letjelica bila jo\u0161 bolje tunirana. S obzirom na upute za tuniranje sljede\u0107i letovi koristili su Stabilised flight mode jer se u njemu trebaju tunirati rate regulatori. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U prvom letu (log_7) tunirali su se svi parametri vezani uz pitch regulatore jer se smatralo da je pitch najkriti\u010dniji od svih drugih pomaka. Nakon \u0161to se relativno istunirao pitch do\u0161lo je do zaklju\u010dka kako su reakcije letjelice preagresivne, odnosno da previ\u0161e agresivno reagira na promjene reference i poreme\u0107aje.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"U drugom i tre\u0107em letu (log_8, log_9) mijenjali su se parametri zasi\u0107enja brzine akceleracije i jerk koji onemogu\u0107avaju letjelici da preagresivno reagira na promjene refenrence ili na poreme\u0107aj. Ovi parametri mijenjali se u Position mode jer je taj mode najvi\u0161e 'automatksi' od dva kori\u0161tena manualn modea \u010dime smo htjeli simularati slu\u010daj kada \u0107e VTOL letjeli automatski tj bez operatora. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U \u010detvrtom i petom letu (log_11 i log_13) ponovo su se i\u0161li tunirati parametri regulatora, ali sada uz nova zasi\u0107enja tj limite na brzine, akceleracije rotacije. Ponovo se krenulo tunirati iz stabilised modea kako je navedeno u uputama. Krenulo se sa tuniranjem po pitch rateu (MC_PITCHRATE_K, MC_PITCHRATE_I, MC_PITCHRATE_D) kako bi dobili \u017eeljeno pona\u0161anje odziva promjene pitcha. Nakon \u0161to se dobilo pribli\u017eno \u017eeljeno pona\u0161anje po pitch rateu, krenulo se na tuniranje roll ratea. Tuniranje roll ratea nije predstavljalo izazov kao \u0161to je bilo tuniranje po pitchu.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"Na kraju dana skinuti su svi konfigurirani paramtri s letjelice i stavljeni na share folder kako bismo imali parametre koji su bolji od prethodnog dana za back up.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"Dvanaesti leta\u010dki dan 11.4.2022\",\n    \"list\": true\n  },\n  {\n    \"text\": \"Nastavak tuniranja regulatora.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"U prvom letu (log_14) krenulo se sa tuniranjem yaw komponenete regulatora odnosno MC_YAW_P. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U drugom letu (log_15) uklju\u010dila se opcija okretanja u vjetar (WV_EN := 1) te se testirala njena funkcionalnost. Uz ovu funkcionalnost PX4 estimira smjer vjetra prema tome u koju stranu se rotira letjelica uspore\u0111uju\u0107i brzinu vrtnju pojedinog od 4 MC motora. Odnosno ako se \u0161alju komande da MC stoji na mjestu tj ne rotira, a MC se zarotira u lijevu stranu zna\u010di da vjetar dolazi s bo\u010dne strane, stoga regulator ispravlja letjelicu prema desno (daje komandu za yaw udesno sve dok vi\u0161e ne bude poreme\u0107aja od vjetra). MC (tj VTOL) je usmjeren u vjetar kada gleda nosom u vjetar jer tada letjelica ne osje\u0107a bo\u010dni poreme\u0107aj. Ova funkcionalnost \u0107e uklju\u010diti na kraju u finalnu letjelicu, no prilikom testiranja nem nije od pomo\u0107i jer regulator mijenja manualne inpute (\u0161to se u konkretnom trenutku ne zna) \u0161to nam mo\u017ee dati krivu interpretaciju rezultata. Uo\u010deno je da funkcionalnost radi kako se od nje i o\u010dekuje.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"U tre\u0107em letu (log_16)  testiralo se u position modeu. Testiralo se kretanje VTOLa u MC fazi lijevo i desno odnosno naprijed i nazad kako bi se snimilo \u0161to vi\u0161e prijelaznih pojava koje bi se mogle analizirati pomo\u0107u PID analize. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"Nakon prou\u010davanja prijelaznih pojava u tre\u0107em letu, u \u010detvrtom letu (log_17) krenulo se na jo\u0161 finije tuniranje regulatora. Tunirali su se roll i pitch rate D komponenete jer se smatralo da D \u010dlan uzrokuje previ\u0161e vibracija u sustavu. Tako\u0111er isklju\u010deno je okretanje u vjetar radi lak\u0161e analize. \",\n    \"list\": false\n  },\n  {\n    \"text\": \"U petom letu (log_22) samo se hoveralo da se prou\u010di utjecaj novih parametara regulatora tokom jednog leta.\",\n    \"list\": false\n  },\n  {\n    \"text\": \"Trinaesti dan 12.4.2022.\",\n    \"list\": true\n  }\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Deseti leta\\u010dki dan 6.4.2022.\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph 'Deseti letacki dan 6.4.2022.' not found.\");\n}\n\nlet current = anchor;\nfor (const para of newParagraphs) {\n  current = current.insertParagraph(para.text, Word.InsertLocation.after);\n  if (para.list) {\n    // Match the existing \"leta\u010dki dan\" heading list style (ListParagraph, numId 3).\n    current.style = \"List Paragraph\";\n    current.attachToList(3, 0);\n  } else {\n    current.style = \"Normal\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Insert the new \"letacki dan\" (flight day) log entries for days 11-13\n# right after the existing \"Deseti letacki dan 6.4.2022.\" paragraph.\n$d = $word.ActiveDocument\n\n$anchor = $null\n$pCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $pCount; $i++) {\n    $cand = $d.Paragraphs.Item($i)\n    if ($cand.Range.Text -like \"*Deseti leta\u010dki dan 6.4.2022.*\") {\n        $anchor = $cand\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph 'Deseti letacki dan 6.4.2022.' not found.\"\n}\n\n# Grab the numbered-list template already used by the \"letacki dan\" headings (numId 3)\n# so the new headings below (Jedanaesti/Dvanaesti/Trinaesti) continue the same list.\n$listTemplate = $anchor.Range.ListFormat.ListTemplate\n\n$current = $anchor\n# paragraph 1\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U ovom danu nastavljalo se tunirati regulatore prema navedenim uputama. Regulatori su se tunirali u Postion modeu. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 2\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U prvom letu (log_0) VTOL letjelica je u\u0161la u nestabilno pona\u0161anje zbog promjene parametara regulatora prije po\u010detka leta. Zaklju\u010deno je da se to dogodilo zbog promjene poja\u010danja na proporcionalnom dijelu pitch rate regulatora te je on nakon leta vra\u0107en na prija\u0161nju vrijednost. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 3\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U drugom letu (log_2) provjeralo se pona\u0161anje letjelice sa staarim parametrima regulatora te je let na kraju bio stabilan \u0161to se i o\u010dekivalo. U ovo letu prou\u010deni su utjecaji promjene parametara vezani uz orijentaciju plo\u010dice, no zaklju\u010deno je kako se ti parametri automatski odrede prilikom kalibracije pa ne postoji potreba za njihovim mijenjanjem. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 4\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U tre\u0107em letu (log_5) nastavilo se s tuniranjem regulatora. Tunirali su se parametri vezani uz pitch i pitchrate odnosno: MC_PITCHRATE_K, MC_PITCH_RATE_D, MC_PITCHRATE_I te MC_PITCH_P. Oni su se mijenjali sve dok se nije dobilo zadovoljavaju\u0107e pona\u0161anje po pitch odnosno dok letjelica nije morala sletiti zbog postotka baterije. Tako\u0111er ovaj dan su spremljeni parametri i stavljeni u share folder gdje se nalaze i logovi. Ovi parametri su karakterizirani kao 'OK parametri' te su se koristili u slu\u010daju kad bi se prilikom tuniranja do\u0161lo do kombinacije parametara pri kojoj bi pona\u0161anje bilo lo\u0161e ili nestabilno.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 5\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Jedanaesti leta\u010dki dan 7.4.2022.\"\n$current.Range.Style = \"List Paragraph\"\n$current.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)\n\n# paragraph 6\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Nastavilo se sa tuniranjem regulatora kako bi se dobili jo\u0161 bolji parametri odnosno kako bi letjelica bila jo\u0161 bolje tunirana. S obzirom na upute za tuniranje sljede\u0107i letovi koristili su Stabilised flight mode jer se u njemu trebaju tunirati rate regulatori. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 7\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U prvom letu (log_7) tunirali su se svi parametri vezani uz pitch regulatore jer se smatralo da je pitch najkriti\u010dniji od svih drugih pomaka. Nakon \u0161to se relativno istunirao pitch do\u0161lo je do zaklju\u010dka kako su reakcije letjelice preagresivne, odnosno da previ\u0161e agresivno reagira na promjene reference i poreme\u0107aje.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 8\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U drugom i tre\u0107em letu (log_8, log_9) mijenjali su se parametri zasi\u0107enja brzine akceleracije i jerk koji onemogu\u0107avaju letjelici da preagresivno reagira na promjene refenrence ili na poreme\u0107aj. Ovi parametri mijenjali se u Position mode jer je taj mode najvi\u0161e 'automatksi' od dva kori\u0161tena manualn modea \u010dime smo htjeli simularati slu\u010daj kada \u0107e VTOL letjeli automatski tj bez operatora. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 9\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U \u010detvrtom i petom letu (log_11 i log_13) ponovo su se i\u0161li tunirati parametri regulatora, ali sada uz nova zasi\u0107enja tj limite na brzine, akceleracije rotacije. Ponovo se krenulo tunirati iz stabilised modea kako je navedeno u uputama. Krenulo se sa tuniranjem po pitch rateu (MC_PITCHRATE_K, MC_PITCHRATE_I, MC_PITCHRATE_D) kako bi dobili \u017eeljeno pona\u0161anje odziva promjene pitcha. Nakon \u0161to se dobilo pribli\u017eno \u017eeljeno pona\u0161anje po pitch rateu, krenulo se na tuniranje roll ratea. Tuniranje roll ratea nije predstavljalo izazov kao \u0161to je bilo tuniranje po pitchu.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 10\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Na kraju dana skinuti su svi konfigurirani paramtri s letjelice i stavljeni na share folder kako bismo imali parametre koji su bolji od prethodnog dana za back up.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 11\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Dvanaesti leta\u010dki dan 11.4.2022\"\n$current.Range.Style = \"List Paragraph\"\n$current.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)\n\n# paragraph 12\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Nastavak tuniranja regulatora.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 13\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U prvom letu (log_14) krenulo se sa tuniranjem yaw komponenete regulatora odnosno MC_YAW_P. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 14\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U drugom letu (log_15) uklju\u010dila se opcija okretanja u vjetar (WV_EN := 1) te se testirala njena funkcionalnost. Uz ovu funkcionalnost PX4 estimira smjer vjetra prema tome u koju stranu se rotira letjelica uspore\u0111uju\u0107i brzinu vrtnju pojedinog od 4 MC motora. Odnosno ako se \u0161alju komande da MC stoji na mjestu tj ne rotira, a MC se zarotira u lijevu stranu zna\u010di da vjetar dolazi s bo\u010dne strane, stoga regulator ispravlja letjelicu prema desno (daje komandu za yaw udesno sve dok vi\u0161e ne bude poreme\u0107aja od vjetra). MC (tj VTOL) je usmjeren u vjetar kada gleda nosom u vjetar jer tada letjelica ne osje\u0107a bo\u010dni poreme\u0107aj. Ova funkcionalnost \u0107e uklju\u010diti na kraju u finalnu letjelicu, no prilikom testiranja nem nije od pomo\u0107i jer regulator mijenja manualne inpute (\u0161to se u konkretnom trenutku ne zna) \u0161to nam mo\u017ee dati krivu interpretaciju rezultata. Uo\u010deno je da funkcionalnost radi kako se od nje i o\u010dekuje.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 15\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U tre\u0107em letu (log_16)  testiralo se u position modeu. Testiralo se kretanje VTOLa u MC fazi lijevo i desno odnosno naprijed i nazad kako bi se snimilo \u0161to vi\u0161e prijelaznih pojava koje bi se mogle analizirati pomo\u0107u PID analize. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 16\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Nakon prou\u010davanja prijelaznih pojava u tre\u0107em letu, u \u010detvrtom letu (log_17) krenulo se na jo\u0161 finije tuniranje regulatora. Tunirali su se roll i pitch rate D komponenete jer se smatralo da D \u010dlan uzrokuje previ\u0161e vibracija u sustavu. Tako\u0111er isklju\u010deno je okretanje u vjetar radi lak\u0161e analize. \"\n$current.Range.Style = \"Normal\"\n\n# paragraph 17\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"U petom letu (log_22) samo se hoveralo da se prou\u010di utjecaj novih parametara regulatora tokom jednog leta.\"\n$current.Range.Style = \"Normal\"\n\n# paragraph 18\n$current.Range.InsertParagraphAfter()\n$current = $current.Next()\n$current.Range.Text = \"Trinaesti dan 12.4.2022.\"\n$current.Range.Style = \"List Paragraph\"\n$current.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)\n\n"}
